$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: bullet list items "Entidade 01" / "Entidade 02" / "Entidade 03"
#         -> "ENTIDADE 01" / "ENTIDADE 02" / "ENTIDADE 03"
# These runs are already bold; only the displayed text changes. Scope each
# replace to its own paragraph so the unrelated "Entidade 01/02/03" text that
# appears later in the document (the sentence with the bracketed list and the
# "Figura N ... Entidade 0N" captions) is left untouched.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Entidade 01" -or $txt -eq "Entidade 02" -or $txt -eq "Entidade 03") {
        $upper = $txt.ToUpper()
        $rng = $d.Range($p.Range.Start, $p.Range.End)
        $rng.Find.ClearFormatting()
        [void]$rng.Find.Execute($txt, $true, $false, $false, $false, $false, $true, 1, $false, $upper, 2)
    }
}

# ---------------------------------------------------------------------------
# Part 2: "Trecho(s): *Nome da caixa* - Entidade 0N" paragraph
#   - "*Nome da caixa* - Entidade 0N" -> "*NOME DA CAIXA* - ENTIDADE 0N"
#     (three separate runs, one per entity)
#   - every run after the "Trecho(s): " label (the three entity runs AND the
#     two "; " separator runs) goes from not-bold to bold
#
# NOTE: this engine coalesces adjacent runs that end up with identical
# formatting whenever a Range.Text/Find-replace text edit touches them. To
# keep the five trailing runs distinct (as the target OOXML requires) we
# temporarily give the three entity runs a one-off Italic flag before editing
# their text - that keeps them from merging with the neighbouring "; " runs,
# which still read Italic=False at that point. Only after all the text edits
# are done do we flip Bold on (and Italic back off) for the whole tail in a
# single property-only pass, which this engine does not use to coalesce runs.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $full = $p.Range
    $t = $full.Text
    if ($t.StartsWith("Trecho(s): ")) {

        $labels = @("*Nome da caixa* - Entidade 01", "*Nome da caixa* - Entidade 02", "*Nome da caixa* - Entidade 03")

        # Step 1: mark the three entity runs Italic so each is distinguishable
        # from its not-yet-edited neighbours during the text edit pass below.
        foreach ($lbl in $labels) {
            $r = $d.Range($full.Start, $full.End)
            $r.Find.ClearFormatting()
            [void]$r.Find.Execute($lbl, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
            $r.Italic = 1
        }

        # Step 2: upper-case the text of each of those three runs.
        foreach ($lbl in $labels) {
            $upperLbl = $lbl.ToUpper()
            $r = $d.Range($full.Start, $full.End)
            $r.Find.ClearFormatting()
            [void]$r.Find.Execute($lbl, $false, $true, $false, $false, $false, $true, 1, $false, "", 0)
            $r.Text = $upperLbl
        }

        # Step 3: bold everything after the "Trecho(s): " label, and clear the
        # temporary Italic marker - a pure property change, so the five
        # distinct runs (3 entity runs + 2 "; " separators) survive intact.
        $lblRng = $d.Range($full.Start, $full.End)
        $lblRng.Find.ClearFormatting()
        [void]$lblRng.Find.Execute("Trecho(s): ", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

        $rest = $d.Range($lblRng.End, $full.End - 1)
        $rest.Bold = 1
        $rest.Italic = 0

        break
    }
}
